# Add two new columns, I (I0) and J (IF), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in row 1, matching the formatting of the existing header cells
# (copy the format from H1, which already carries the header style).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J16.
$values = @(
    @(1, 5),
    @(1, 6),
    @(1, 4),
    @(1, 6),
    @(1, 4),
    @(1, 4),
    @(1, 7),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(1, 3),
    @(1, 4),
    @(1, 5),
    @(2, 6),
    @(5, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
